$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 6's table: switch its table style to the built-in style
#    {6462257C-74E7-45B3-8B3C-67C4FA5EC96B}
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6462257C-74E7-45B3-8B3C-67C4FA5EC96B}", $true)
    }
}

# ---------------------------------------------------------------------------
# 2) Deck theme (ppt/theme/theme1.xml, used by the slide master): swap the
#    "Integral" colour scheme for the stock "Office" colour scheme.
#    The font scheme / format scheme are already identical between the two
#    themes, so only the 12 theme colours need to move.
#    PowerPoint's RGB long uses 0x00BBGGRR ordering.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
